$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rename mapping in column A: the prefix value for row 2 changed,
# and a new prefix was added for row 3 (rename logic now covers another group).
$ws.Range("A2").Value = "L8258_T11"
$ws.Range("A3").Value = "C8161_X11"
